# Auto-generated Excel COM-interop script
# Applies numeric value updates (and a few cell add/remove ops) to the
# 'Masamune_Profits' market-data columns (H-N) across all 8 sheets,
# matching the scheduled-runner refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H64").Value = 23028.611
$ws.Range("I64").Value = 61179.41
$ws.Range("J64").Value = 2761
$ws.Range("K64").Value = 61179.41
$ws.Range("L64").Value = 2761
$ws.Range("M64").Value = -60931.41
$ws.Range("N64").Value = -3257
$ws.Range("H67").Value = 23028.611
$ws.Range("I67").Value = 61179.41
$ws.Range("J67").Value = 2761
$ws.Range("K67").Value = 61179.41
$ws.Range("L67").Value = 2761
$ws.Range("M67").Value = -60321.41
$ws.Range("N67").Value = -4477
$ws.Range("H98").Value = 36830.773
$ws.Range("I98").Value = 1005
$ws.Range("J98").Value = 66685.586
$ws.Range("K98").Value = 1005
$ws.Range("L98").Value = 66685.586
$ws.Range("M98").Value = 493
$ws.Range("N98").Value = -69681.586
$ws.Range("H100").Value = 1393.6154
$ws.Range("I100").Value = 1311
$ws.Range("J100").Value = 1579.5
$ws.Range("K100").Value = 1311
$ws.Range("L100").Value = 1579.5
$ws.Range("M100").Value = -770
$ws.Range("N100").Value = -2661.5
$ws.Range("H108").Value = 31217.6
$ws.Range("J108").Value = 31217.6
$ws.Range("L108").Value = 31217.6
$ws.Range("N108").Value = -38897.6
$ws.Range("H109").Value = 28919
$ws.Range("J109").Value = 28919
$ws.Range("L109").Value = 28919
$ws.Range("N109").Value = -31693
$ws.Range("H120").Value = 49706
$ws.Range("J120").Value = 49706
$ws.Range("L120").Value = 49706
$ws.Range("N120").Value = -59382
$ws.Range("H122").Value = 36830.773
$ws.Range("I122").Value = 1005
$ws.Range("J122").Value = 66685.586
$ws.Range("K122").Value = 3015
$ws.Range("L122").Value = 200056.758
$ws.Range("M122").Value = -565
$ws.Range("N122").Value = -204956.758
$ws.Range("H123").Value = 32564.615
$ws.Range("J123").Value = 32564.615
$ws.Range("L123").Value = 32564.615
$ws.Range("N123").Value = -42364.61500000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 26989.35
$ws.Range("I32").Value = 27328.25
$ws.Range("J32").Value = 24278.143
$ws.Range("K32").Value = 27328.25
$ws.Range("L32").Value = 24278.143
$ws.Range("M32").Value = -27041.25
$ws.Range("N32").Value = -24852.143
$ws.Range("H80").Value = 50615.145
$ws.Range("J80").Value = 50615.145
$ws.Range("L80").Value = 50615.145
$ws.Range("N80").Value = -52611.145
$ws.Range("H83").Value = 50615.145
$ws.Range("J83").Value = 50615.145
$ws.Range("L83").Value = 151845.435
$ws.Range("N83").Value = -161829.435
$ws.Range("H120").Value = 42197.332
$ws.Range("J120").Value = 42197.332
$ws.Range("L120").Value = 42197.332
$ws.Range("N120").Value = -51873.332
$ws.Range("H132").Value = 12821944
$ws.Range("I132").Value = 27778844
$ws.Range("J132").Value = 1743.4286
$ws.Range("K132").Value = 83336532
$ws.Range("L132").Value = 5230.2858
$ws.Range("M132").Value = -83334002
$ws.Range("N132").Value = -10290.2858
$ws.Range("H138").Value = 46088.75
$ws.Range("J138").Value = 46088.75
$ws.Range("L138").Value = 46088.75
$ws.Range("N138").Value = -56368.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H117").Value = 44999
$ws.Range("J117").Value = 44999
$ws.Range("L117").Value = 44999
$ws.Range("N117").Value = -54177
$ws.Range("H119").Value = 47992
$ws.Range("J119").Value = 47992
$ws.Range("L119").Value = 47992
$ws.Range("N119").Value = -57668
$ws.Range("H120").Value = 48753
$ws.Range("J120").Value = 48753
$ws.Range("L120").Value = 48753
$ws.Range("N120").Value = -58429
$ws.Range("H134").Value = 3928.8125
$ws.Range("I134").Value = 3017.111
$ws.Range("J134").Value = 4193.5
$ws.Range("K134").Value = 9051.332999999999
$ws.Range("L134").Value = 12580.5
$ws.Range("M134").Value = -6516.332999999999
$ws.Range("N134").Value = -17650.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H104").Value = 30872.875
$ws.Range("J104").Value = 30872.875
$ws.Range("L104").Value = 30872.875
$ws.Range("N104").Value = -36114.875
$ws.Range("H109").Value = 27701.7
$ws.Range("J109").Value = 27701.7
$ws.Range("L109").Value = 27701.7
$ws.Range("N109").Value = -29781.7
$ws.Range("H115").Value = 28958.334
$ws.Range("J115").Value = 28958.334
$ws.Range("L115").Value = 28958.334
$ws.Range("N115").Value = -31308.334
$ws.Range("H116").Value = 42362.25
$ws.Range("J116").Value = 42362.25
$ws.Range("L116").Value = 42362.25
$ws.Range("N116").Value = -51540.25
$ws.Range("H120").Value = 33578.7
$ws.Range("J120").Value = 33578.7
$ws.Range("L120").Value = 33578.7
$ws.Range("N120").Value = -40836.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 3814.7837
$ws.Range("I5").Value = 12860.875
$ws.Range("J5").Value = 1319.3103
$ws.Range("K5").Value = 38582.625
$ws.Range("L5").Value = 3957.9309
$ws.Range("M5").Value = -38470.625
$ws.Range("N5").Value = -4181.9309
$ws.Range("H94").Value = 5250
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 5250
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 15750
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -17102
$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -19118
$ws.Range("H98").Value = 806.0833
$ws.Range("I98").Value = 528.1667
$ws.Range("J98").Value = 1084
$ws.Range("K98").Value = 1584.5001
$ws.Range("L98").Value = 3252
$ws.Range("M98").Value = -86.50009999999997
$ws.Range("N98").Value = -6248
$ws.Range("H125").Value = 4345.3335
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4345.3335
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 13036.0005
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -22876.0005
$ws.Range("H126").Value = 24275.268
$ws.Range("I126").Value = 75852.5
$ws.Range("J126").Value = 5519.909
$ws.Range("K126").Value = 227557.5
$ws.Range("L126").Value = 16559.727
$ws.Range("M126").Value = -222617.5
$ws.Range("N126").Value = -26439.727
$ws.Range("H133").Value = 7397.0835
$ws.Range("I133").Value = 8425
$ws.Range("J133").Value = 6883.125
$ws.Range("K133").Value = 25275
$ws.Range("L133").Value = 20649.375
$ws.Range("M133").Value = -20215
$ws.Range("N133").Value = -30769.375
$ws.Range("H134").Value = 47624820
$ws.Range("I134").Value = 100004450
$ws.Range("J134").Value = 6982.273
$ws.Range("K134").Value = 300013350
$ws.Range("L134").Value = 20946.819
$ws.Range("M134").Value = -300008280
$ws.Range("N134").Value = -31086.819
$ws.Range("H135").Value = 3814.7837
$ws.Range("I135").Value = 12860.875
$ws.Range("J135").Value = 1319.3103
$ws.Range("K135").Value = 115747.875
$ws.Range("L135").Value = 11873.7927
$ws.Range("M135").Value = -113212.875
$ws.Range("N135").Value = -16943.7927
$ws.Range("H139").Value = 152272.9
$ws.Range("I139").Value = 216244.28
$ws.Range("J139").Value = 3006.3333
$ws.Range("K139").Value = 648732.84
$ws.Range("L139").Value = 9018.999899999999
$ws.Range("M139").Value = -643592.84
$ws.Range("N139").Value = -19298.9999
$ws.Range("H140").Value = 2181.4
$ws.Range("I140").Value = 1757.4166
$ws.Range("J140").Value = 3877.3333
$ws.Range("K140").Value = 5272.2498
$ws.Range("L140").Value = 11631.9999
$ws.Range("M140").Value = -92.2497999999996
$ws.Range("N140").Value = -21991.9999
$ws.Range("H141").Value = 83337630
$ws.Range("I141").Value = 111115176
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 333345528
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -333340348
$ws.Range("N141").Value = -25360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H104").Value = 42740.4
$ws.Range("J104").Value = 42740.4
$ws.Range("L104").Value = 42740.4
$ws.Range("N104").Value = -49728.4
$ws.Range("H110").Value = 47694
$ws.Range("J110").Value = 47694
$ws.Range("L110").Value = 47694
$ws.Range("N110").Value = -55874
$ws.Range("H118").Value = 34815
$ws.Range("J118").Value = 34815
$ws.Range("L118").Value = 34815
$ws.Range("N118").Value = -38129
$ws.Range("H130").Value = 45985.25
$ws.Range("J130").Value = 45985.25
$ws.Range("L130").Value = 45985.25
$ws.Range("N130").Value = -56025.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H111").Value = 43848
$ws.Range("J111").Value = 43848
$ws.Range("L111").Value = 43848
$ws.Range("N111").Value = -52028
$ws.Range("H121").Value = 40277.332
$ws.Range("J121").Value = 40277.332
$ws.Range("L121").Value = 40277.332
$ws.Range("N121").Value = -43771.332

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H16").Value = 40718.5
$ws.Range("J16").Value = 40718.5
$ws.Range("L16").Value = 40718.5
$ws.Range("N16").Value = -41302.5
$ws.Range("H119").Value = 48690
$ws.Range("J119").Value = 48690
$ws.Range("L119").Value = 48690
$ws.Range("N119").Value = -58366
$ws.Range("H136").Value = 18668.932
$ws.Range("I136").Value = 44243.61
$ws.Range("J136").Value = 1862.7142
$ws.Range("K136").Value = 132730.83
$ws.Range("L136").Value = 5588.142599999999
$ws.Range("M136").Value = -130180.83
$ws.Range("N136").Value = -10688.1426
$ws.Range("H137").Value = 49333
$ws.Range("J137").Value = 49333
$ws.Range("L137").Value = 49333
$ws.Range("N137").Value = -59533

Write-Output "Applied updates across 8 sheets (240 set, 2 added, 2 cleared)."